# Update the NATMI LR-pairs sheet (Col1a1-Gp6) with recomputed TPM-based
# values. The sending/target cluster table grows from 4 rows (FAPs only as
# target) to 8 rows (FAPs and MuSCs as target for every sending cluster),
# and the numeric statistics (G:T) are recalculated throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A..T as in the header row.
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @{
    2 = @("ECs","Col1a1","Gp6","FAPs",3,1,26.81310566666667,80.439317,0.004518206005002021,0.004518206005002021,2,0.6666666666666666,0.053317,0.159951,0.5209400637695697,0.5209400637695697,1.429594354829667,12.866349193467,0.002353714524369806,0.002353714524369806)
    3 = @("ECs","Col1a1","Gp6","MuSCs",3,1,26.81310566666667,80.439317,0.004518206005002021,0.004518206005002021,1,0.3333333333333333,0.04903066666666667,0.147092,0.4790599362304302,0.4790599362304302,1.314664446240445,11.831980016164,0.002164491480632215,0.002164491480632215)
    4 = @("FAPs","Col1a1","Gp6","FAPs",3,1,5771.873535333333,17315.620606,0.9726032482643521,0.9726032482643523,2,0.6666666666666666,0.053317,0.159951,0.5209400637695697,0.5209400637695697,307.7389812833674,2769.650831550306,0.5066679981733222,0.5066679981733223)
    5 = @("FAPs","Col1a1","Gp6","MuSCs",3,1,5771.873535333333,17315.620606,0.9726032482643521,0.9726032482643523,1,0.3333333333333333,0.04903066666666667,0.147092,0.4790599362304302,0.4790599362304302,282.9988073530835,2546.989266177752,0.4659352500910298,0.4659352500910299)
    6 = @("MuSCs","Col1a1","Gp6","FAPs",3,1,132.4457753333333,397.337326,0.02231808970163987,0.02231808970163988,2,0.6666666666666666,0.053317,0.159951,0.5209400637695697,0.5209400637695697,7.061611403447333,63.554502631026,0.01162638707238725,0.01162638707238725)
    7 = @("MuSCs","Col1a1","Gp6","MuSCs",3,1,132.4457753333333,397.337326,0.02231808970163987,0.02231808970163988,1,0.3333333333333333,0.04903066666666667,0.147092,0.4790599362304302,0.4790599362304302,6.493904661776889,58.445141955992,0.01069170262925262,0.01069170262925262)
    8 = @("Resolving-Mac","Col1a1","Gp6","FAPs",3,1,3.326003,9.978009,0.0005604560290058679,0.000560456029005868,2,0.6666666666666666,0.053317,0.159951,0.5209400637695697,0.5209400637695697,0.177332501951,1.595992517559,0.0002919639994903566,0.0002919639994903567)
    9 = @("Resolving-Mac","Col1a1","Gp6","MuSCs",3,1,3.326003,9.978009,0.0005604560290058679,0.000560456029005868,1,0.3333333333333333,0.04903066666666667,0.147092,0.4790599362304302,0.4790599362304302,0.1630761444253333,1.467685299828,0.0002684920295155112,0.0002684920295155113)
}

foreach ($r in 2..9) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$r").Value = $values[$i]
    }
}
